# ---------------------------------------------------------------------------
# Applies two changes to the deck:
#
# 1. Slide 16's table (the only table in the deck) gets a new table style
#    (tableStyleId {383A89FA-A425-4FEA-8DCC-2219C790C396} ->
#     {A89D8E3B-3947-4943-B46E-B43CF6994D1B}).
#
# 2. The presentation's theme colour palette is swapped from the custom
#    "Integral" palette over to the stock Office palette (the "Integral"
#    and "Office Theme" theme parts trade places / the deck's clrScheme is
#    repointed at the default Office colours).
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Table style -----------------------------------------------------
$s = $p.Slides.Item(16)
$tableShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasTable) {
        $tableShape = $candidate
        break
    }
}
$tbl = $tableShape.Table
$tbl.ApplyStyle("{A89D8E3B-3947-4943-B46E-B43CF6994D1B}")

# --- 2. Theme colours: Integral -> Office --------------------------------
# Colour order matches the 12-slot theme colour scheme:
#   dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
# RGB values are packed as 0xBBGGRR (OLE COLORREF), matching the
# PowerPoint COM convention for the .RGB property.
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$slide1 = $p.Slides.Item(1)
$themeColors = $slide1.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officeColors[$i - 1]
}
